# "Generate Report for Handback"
#
# The 05997961-... file's handback transform failed because the handback
# filename (bj1kiggw.jey) didn't match the handoff filename. Reflect that
# in the localization-status report:
#   - Status for that file flips from "Ready for handoff" to
#     "Handback transform failed" (on every sheet that shows it).
#   - Each language sheet gets an "Error Detail" explaining the mismatch.
#   - The "Error Detail" column is widened so the message is readable.

$wb = $excel.ActiveWorkbook

$failedStatus = "Handback transform failed"

$zhError = "Handback file name: bj1kiggw.jey is different with handoff file name: 05997961-87b0-464c-9fca-8982507ba71e.eae5ceb3214052adef449791dc64ea8f546253b8.zh-cn."
$deError = "Handback file name: bj1kiggw.jey is different with handoff file name: 05997961-87b0-464c-9fca-8982507ba71e.eae5ceb3214052adef449791dc64ea8f546253b8.de-de."

# --- Overview sheet: status shown for the 05997961 row (zh-cn / de-de cols) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $failedStatus
$overview.Range("F3").Value = $failedStatus

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $failedStatus
$zh.Range("P3").Value = $zhError
$zh.Columns.Item(16).ColumnWidth = 39.2

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $failedStatus
$de.Range("P3").Value = $deError
$de.Columns.Item(16).ColumnWidth = 39.2
